$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# JobHistory test-suite stats: one fewer "partially automated" case tracked
# in B4 (EmailNotification / JobHistoryErrors row). G5 (SUM of column B) and
# G6 (G5/G4 pass-rate) are formulas and recalc automatically.
$ws.Range("B4").Value = 1

# Move the cursor/selection to B5, matching the saved cell pointer.
$ws.Range("B5").Select()
